$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header row (row 1) contents - cells stay empty, no shifting
$ws.Rows.Item(1).ClearContents()

# Clear the Description column (column E) contents - cells stay empty, no shifting
$ws.Columns.Item(5).ClearContents()

# Update cell values for the data rows (rows 2-6, columns A-D).
# Values are forced to be stored as text (e.g. "25" not 25) by briefly
# applying a text number format and clearing formats back afterward so
# no visible style is left on the cells.
$data = $ws.Range("A2:D6")
$data.NumberFormat = "@"

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "cita"
$ws.Range("C2").Value = "25"
$ws.Range("D2").Value = "25"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "cita"
$ws.Range("C3").Value = "vh"
$ws.Range("D3").Value = "hv"

$ws.Range("A4").Value = "1"
$ws.Range("B4").Value = "para"
$ws.Range("C4").Value = "453"
$ws.Range("D4").Value = "42"

$ws.Range("A5").Value = "1"
$ws.Range("B5").Value = "hy"
$ws.Range("C5").Value = "45"
$ws.Range("D5").Value = "54"

$ws.Range("A6").Value = "1"
$ws.Range("B6").Value = "cita"
$ws.Range("C6").Value = "nk"
$ws.Range("D6").Value = "njm"

$data.ClearFormats()
